# Fix gap filling indexing bug:
# Add the missing "CH4_flux" value to the Vars_to_fill column (C) on row 4
# of the Berge_MDS sheet, and move the active selection to C5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Berge_MDS")

$ws.Range("C4").Value = "CH4_flux"

$ws.Activate()
$ws.Range("C5").Select()
